$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the three daily values that were previously placeholders (0)
$ws.Range("B57").Value = 5
$ws.Range("B58").Value = 9
$ws.Range("B59").Value = 1

# Match the updated view state: scroll position + current selection
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E61").Select() | Out-Null
